$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 12: date label "2025-08-30" as literal text (same shape as the other date cells) ---
$ws.Range("A12").Formula = '="2025-08-30"'
$ws.Range("A12").Copy()
$ws.Range("A12").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Revised nowcast figures for rows 2-11, plus the new row 12 ---
# Row 2
$ws.Range("B2").Value = 0.29582471055420528
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

# Row 3
$ws.Range("B3").Value = 0.29611163418938269
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.00094385170027047464
$ws.Range("E3").Value = 0.000012604605672936883
$ws.Range("F3").Value = -0.000026605782403536365
$ws.Range("G3").Value = -0.00016735052225396214
$ws.Range("H3").Value = -0.0000010299678870835377
$ws.Range("I3").Value = -0.00058340629621847116
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.00010885989799702678

# Row 4
$ws.Range("B4").Value = 0.28373080472794981
$ws.Range("C4").Value = -0.011085349261032563
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = -0.00003275447347262766
$ws.Range("F4").Value = 0.00000050100608448572638
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = -0.000010834873242485248
$ws.Range("I4").Value = -0.0020566461662239431
$ws.Range("J4").Value = 0.00085480780886082268
$ws.Range("K4").Value = -0.000050553502406569528

# Row 5
$ws.Range("B5").Value = 0.28947436878975819
$ws.Range("C5").Value = 0.010254636325539686
$ws.Range("D5").Value = 0.0013323996584479374
$ws.Range("E5").Value = -0.00046684675408094583
$ws.Range("F5").Value = -0.00027905171241852993
$ws.Range("G5").Value = -0.0046758603353293923
$ws.Range("H5").Value = -0.00010637500542164646
$ws.Range("I5").Value = -0.000089292469806835401
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -0.00022604564512190839

# Row 6
$ws.Range("B6").Value = 0.47736865043376087
$ws.Range("C6").Value = 0.21508514223227906
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.00060460097928950426
$ws.Range("F6").Value = 0.000032128137756242341
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -0.00037045638209440675
$ws.Range("I6").Value = -0.031228579107887695
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.003771445784659988

# Row 7
$ws.Range("B7").Value = 0.46521019738093039
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.030745811586871844
$ws.Range("E7").Value = -0.001270196775723821
$ws.Range("F7").Value = -0.0050944160501147709
$ws.Range("G7").Value = 0.023363815847833643
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.0019938287587492353
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -0.00040567324670293115

# Row 8
$ws.Range("B8").Value = 0.13507050165897355
$ws.Range("C8").Value = -0.33065893061262558
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = -0.00040562550247949178
$ws.Range("F8").Value = -0.0067309152813595022
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0.00038943555045950473
$ws.Range("I8").Value = 0.004790199569762954
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0.0024761405542852444

# Row 9
$ws.Range("B9").Value = -0.099298831436153312
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = -0.037204726125556538
$ws.Range("E9").Value = -0.027967263323098592
$ws.Range("F9").Value = -0.17375710225479618
$ws.Range("G9").Value = 0.0076654367780558431
$ws.Range("H9").Value = -0.0057040009342530732
$ws.Range("I9").Value = 0.0022523177579940067
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0.00034600500652764987

# Row 10
$ws.Range("B10").Value = 0.20945158603866781
$ws.Range("C10").Value = 0.30724144800612396
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = -0.0013429917906939192
$ws.Range("F10").Value = 0.0011056104287888589
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = -0.0011002900973423857
$ws.Range("I10").Value = 0.025029839058208436
$ws.Range("J10").Value = -0.019005804488106137
$ws.Range("K10").Value = -0.0031773936421577442

# Row 11
$ws.Range("B11").Value = 0.44004877102480283
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0.050749500380081546
$ws.Range("E11").Value = 0.024037652493940492
$ws.Range("F11").Value = 0.17907134085366511
$ws.Range("G11").Value = 0.0099036637743366435
$ws.Range("H11").Value = 0.00055865197186969872
$ws.Range("I11").Value = 0.015016943068166914
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = -0.048740567555925374

# Row 12
$ws.Range("B12").Value = 0.43186519303654303
$ws.Range("C12").Value = 0.072200316379180479
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = -0.00048622990055180828
$ws.Range("F12").Value = 0.0002491055089423349
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = -0.00020677182640028372
$ws.Range("I12").Value = -0.060604038008273187
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = -0.019335960141157349

# --- Column width tweaks (col B/C) ---
$ws.Columns.Item(2).ColumnWidth = 14.24609375
$ws.Columns.Item(3).ColumnWidth = 14.24609375
